# GDP_Steo.xlsx update — "Updated data feeds except mine netback and interim"
#
# The STEO GDP feed was refreshed: the most recent quarters (2022 Q3/Q4 and
# all of 2023) received new GDP values, and the cells that used to carry the
# "estimate" (italic) formatting for the 2022-Q4 through 2023 rows are now
# treated as normal/actual data, so their formatting is switched from the
# italic style to the same bold, right-aligned style already used by the
# rest of column C (as seen on e.g. C176).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) Re-format C177:C193 so they match the (non-italic) style already used
#    by C122:C176, instead of the italic "estimate" style they had before.
# ---------------------------------------------------------------------
$ws.Range("C176").Copy()
$ws.Range("C177:C193").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Write the refreshed GDP values for 2022Q3 through 2023Q4 (rows 173-193)
# ---------------------------------------------------------------------
$newValues = @{
    173 = 19681.682000000001
    174 = 19681.682000000001
    175 = 19681.682000000001
    176 = 19694.757407000001
    177 = 19707.065852
    178 = 19722.84
    179 = 19741.84
    180 = 19764.71
    181 = 19791.21
    182 = 19822.919999999998
    183 = 19855.52
    184 = 19890.580000000002
    185 = 19930.12
    186 = 19968.59
    187 = 20007.990000000002
    188 = 20051.32
    189 = 20090.37
    190 = 20128.13
    191 = 20166.28
    192 = 20200.169999999998
    193 = 20231.509999999998
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value2 = $newValues[$row]
}

# ---------------------------------------------------------------------
# 3) Move the selection/scroll position to where the author left off
#    (scrolled down near the bottom of the table, cell D176 selected).
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 178
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D176").Select()
